$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '61.517.02'
$ws.Range('E2').Value = '  -1.90%  '
$ws.Range('D3').Value = '2.994.07'
$ws.Range('E3').Value = '  -1.08%  '
$ws.Range('E4').Value = '  +0.01%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '595.47'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  +2.15%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '143.59'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  -3.53%  '
$ws.Range('E7').Value = '  +0.05%  '
$ws.Range('E8').Value = '  -0.76%  '
$ws.Range('D9').Value = '2.987.45'
$ws.Range('E9').Value = '  -1.32%  '
$ws.Range('E10').Value = '  -2.68%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '5.91'
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '  +3.77%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.462'
$ws.Range('D12').ClearFormats()
$ws.Range('E12').Value = '  +4.24%  '
$ws.Range('E13').Value = '  -1.25%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '34.22'
$ws.Range('D14').ClearFormats()
$ws.Range('E15').Value = '  +2.11%  '
$ws.Range('D16').Value = '3.485.47'
$ws.Range('E16').Value = '  -1.16%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '7.00'
$ws.Range('D17').ClearFormats()
$ws.Range('E17').Value = '  -1.05%  '
$ws.Range('D18').Value = '61.494.61'
$ws.Range('E18').Value = '  -1.86%  '
$ws.Range('D19').Value = '2.994.38'
$ws.Range('E19').Value = '  -1.01%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '453.58'
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  -3.10%  '
$ws.Range('E21').Value = '  -0.53%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.687'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  -0.73%  '
$ws.Range('E23').Value = '  -0.78%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '82.27'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  +1.64%  '
$ws.Range('E25').Value = '  -6.90%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '12.08'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  -2.58%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '10.36'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  +0.03%  '
$ws.Range('E28').Value = '  +0.15%  '
$ws.Range('E29').Value = '  +1.93%  '
$ws.Range('E30').Value = '  +0.04%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '6.98'
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  -3.84%  '
$ws.Range('E32').Value = '  -3.51%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '27.60'
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  +0.60%  '
$ws.Range('E34').Value = '  +0.11%  '
$ws.Range('D35').Value = '0.0₃0814'
$ws.Range('E35').Value = '  +2.41%  '
$ws.Range('E36').Value = '  -1.47%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '5.75'
$ws.Range('D37').ClearFormats()
$ws.Range('E37').Value = '  -0.62%  '
$ws.Range('E38').Value = '  -3.72%  '
$ws.Range('B39').Value = 'OKB'
$ws.Range('C39').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '50.31'
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  +0.17%  '
$ws.Range('B40').Value = 'Cosmos'
$ws.Range('C40').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '9.17'
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = '  +1.81%  '
$ws.Range('E41').Value = '  +7.88%  '
$ws.Range('E42').Value = '  -3.13%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '397.34'
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  -5.99%  '
$ws.Range('E44').Value = '  -0.99%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '38.64'
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  +2.32%  '
$ws.Range('E46').Value = '  -5.47%  '
$ws.Range('D47').Value = '2.718.04'
$ws.Range('E47').Value = '  -3.03%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '133.63'
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  +3.51%  '
$ws.Range('E50').Value = '  -0.76%  '
$ws.Range('E51').Value = '  +0.92%  '
